# "fixed some errors in results"
# In the ablation-timing / RF-lesion-count columns (AZ:BD), a value of 0
# was being used to mean "not applicable" (procedure not performed on that
# vein / no gap re-touch-up needed, etc.) which is indistinguishable from a
# genuine zero. Replace those literal-zero cells with "/" (not applicable)
# in the affected columns, and also mark a few missing probnp (O) readings
# that were left blank by mistake.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 60
$cols = @("AZ", "BA", "BB", "BC", "BD")

foreach ($col in $cols) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Range("$col$r")
        $v = $cell.Value2
        if ([string]::IsNullOrEmpty($v)) {
            continue
        }
        if ($v -eq 0) {
            $cell.Value = "/"
        }
    }
}

# A couple of genuinely-missing probnp readings that were left blank;
# mark them explicitly as not-available too.
$ws.Range("O10").Value = "/"
$ws.Range("O13").Value = "/"
$ws.Range("O14").Value = "/"

# Two cells in row 57 (ablate_reisolization_time / ablate_removal_time_dormant)
# were left completely empty instead of being marked not-applicable.
$ws.Range("BA57").Value = "/"
$ws.Range("BB57").Value = "/"
